# Update "想去人数" (column F) values across worksheets per gh-pages regeneration
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 1156
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 7141
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 5484
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 6254
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 412
$ws.Range("F18").Value = 31
$ws.Range("F19").Value = 0
$ws.Range("F21").Value = 208
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 1927
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 2122
$ws.Range("F32").Value = 88
$ws.Range("F34").Value = 1034
$ws.Range("F35").Value = 16
$ws.Range("F37").Value = 311
$ws.Range("F38").Value = 62
$ws.Range("F39").Value = 5222
$ws.Range("F41").Value = 666
$ws.Range("F44").Value = 1102
$ws.Range("F45").Value = 1078
$ws.Range("F47").Value = 1374
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 1096

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 32
$ws.Range("F10").Value = 194
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F14").Value = 5
$ws.Range("F16").Value = 1

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 37
$ws.Range("F4").Value = 68
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 16
$ws.Range("F7").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = 185
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 5484
$ws.Range("F18").Value = 6255
$ws.Range("F19").Value = 6255
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 412
$ws.Range("F22").Value = 274
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 194
$ws.Range("F28").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 1034
$ws.Range("F38").Value = 2085
$ws.Range("F40").Value = 5222
$ws.Range("F41").Value = 1202
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 111
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 1374
$ws.Range("F49").Value = 0
$ws.Range("F50").Value = 0
